$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "usuario" column for the anulación (cancellation) rows 5-7
$ws.Range("C5").Value = "apellegrini"
$ws.Range("C6").Value = "apellegrini"
$ws.Range("C7").Value = "apellegrini"

# Update the NroSiniestro values (text cells, preserve exact spacing with a
# leading apostrophe so Excel keeps treating them as literal text)
$ws.Range("F2").Value = "'1220194200684 "
$ws.Range("F3").Value = "'1120194100442 "
$ws.Range("F4").Value = "'0420194406812"
$ws.Range("F5").Value = "'0420172008629  "
$ws.Range("F6").Value = "'1220170301442  "
$ws.Range("F7").Value = "'1120170200942  "

# Update the active selection to match the saved view state
$ws.Range("F10").Select()
